$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts existing D:K to E:L)
$ws.Columns("D:D").Insert()

# Copy formatting (number formats/styles) from column E (the old column D, now
# shifted right) down into the new column D, bounded to the used data rows so
# the sheet dimension / used range isn't blown out to the full column height.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows that have no data in columns D:L at all (section/header labels only)
# shouldn't end up with a stray empty, styled D cell from the paste above.
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# Populate new column D with the new quarter values
$ws.Cells.Item(7, 4).Value = 43373
$ws.Cells.Item(8, 4).Value = 34800
$ws.Cells.Item(9, 4).Value = 12100
$ws.Cells.Item(10, 4).Value = 22700
$ws.Cells.Item(12, 4).Value = 8200
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = -2100
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(17, 4).Value = 45700
$ws.Cells.Item(18, 4).Value = -10900
$ws.Cells.Item(20, 4).Value = 8200
$ws.Cells.Item(21, 4).Value = -700
$ws.Cells.Item(22, 4).Value = 7200
$ws.Cells.Item(23, 4).Value = -10000
$ws.Cells.Item(24, 4).Value = 200
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = -10100
$ws.Cells.Item(27, 4).Value = -10100
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = -8200
$ws.Cells.Item(33, 4).Value = -10100
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = -10100
$ws.Cells.Item(38, 4).Value = 43373
$ws.Cells.Item(41, 4).Value = 38300
$ws.Cells.Item(42, 4).Value = "NA"
$ws.Cells.Item(43, 4).Value = 22800
$ws.Cells.Item(44, 4).Value = 42300
$ws.Cells.Item(45, 4).Value = 6700
$ws.Cells.Item(46, 4).Value = 110100
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 16900
$ws.Cells.Item(49, 4).Value = 198200
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 1100
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 326300
$ws.Cells.Item(57, 4).Value = 11300
$ws.Cells.Item(58, 4).Value = 28000
$ws.Cells.Item(59, 4).Value = 31600
$ws.Cells.Item(60, 4).Value = 71000
$ws.Cells.Item(61, 4).Value = 189100
$ws.Cells.Item(62, 4).Value = 23100
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 283200
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = -573800
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 43100
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43373
$ws.Cells.Item(81, 4).Value = -10100
$ws.Cells.Item(83, 4).Value = 2000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = -6500
$ws.Cells.Item(91, 4).Value = -100
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -100
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = 11300
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(102, 4).Value = 4800
